$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.592524647712708
$ws.Range("B1").Value = 3.671234130859375
$ws.Range("C1").Value = 5.452422618865967
$ws.Range("D1").Value = 1.347572088241577
$ws.Range("E1").Value = 0.78338223695755
